$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (Coin names, URLs, Volume strings) -- safe to assign directly
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("E43").Value = '42CEJICEJIBestin24h'
$ws.Range("E48").Value = '47BOLOBOLO'

# Numeric-looking price strings must be forced to Text so Excel keeps the exact
# literal representation (trailing zeros, precision) instead of parsing as a Double.
# NumberFormat is reset back to General afterwards so no stray cell style is left behind.
$priceCells = @("D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D40", "D41", "D44", "D45", "D46", "D47", "D48")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D3").Value = '24.20'
$ws.Range("D4").Value = '5.294'
$ws.Range("D5").Value = '0.05794'
$ws.Range("D6").Value = '6.478'
$ws.Range("D7").Value = '3.136'
$ws.Range("D8").Value = '0.8168'
$ws.Range("D9").Value = '0.8771'
$ws.Range("D10").Value = '0.1380'
$ws.Range("D11").Value = '0.07004'
$ws.Range("D12").Value = '0.03134'
$ws.Range("D13").Value = '0.02917'
$ws.Range("D14").Value = '0.09415'
$ws.Range("D15").Value = '3.748'
$ws.Range("D16").Value = '0.001529'
$ws.Range("D17").Value = '0.04668'
$ws.Range("D18").Value = '0.0005995'
$ws.Range("D19").Value = '0.006052'
$ws.Range("D20").Value = '0.001245'
$ws.Range("D21").Value = '0.004664'
$ws.Range("D22").Value = '0.00006102'
$ws.Range("D23").Value = '3.533'
$ws.Range("D24").Value = '2.144'
$ws.Range("D25").Value = '0.3189'
$ws.Range("D28").Value = '0.0002334'
$ws.Range("D40").Value = '0.03723'
$ws.Range("D41").Value = '0.006372'
$ws.Range("D44").Value = '0.007744'
$ws.Range("D45").Value = '0.00005275'
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("D47").Value = '0.4034'
$ws.Range("D48").Value = '0.002335'

foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "General"
    $ws.Range($ref).Style = "Normal"
}
